$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: swap the SamplePortion / Result header order (H1/I1)
$ws.Range("H1").Value = "SamplePortion"
$ws.Range("I1").Value = "Result"

# Row 2: add unit to the float type for the swapped SamplePortion/Result columns
$ws.Range("H2").Value = "#float,  unit:mg"
$ws.Range("I2").Value = "#float,  unit:mg"

# Row 3: new French description/enum row
$ws.Range("A3").Value = "#Manipulateur"
$ws.Range("B3").Value = "#Desc:IdentifiantEchantillon"
$ws.Range("C3").Value = "#Date"
$ws.Range("D3").Value = "#ModeOderatoireLaboratoire"
$ws.Range("E3").Value = "#AppareilLogicielCritique"
$ws.Range("F3").Value = "#ProduitCritique"
$ws.Range("G3").Value = "#LieuStockageDonneesBrutes"
$ws.Range("H3").Value = "#PriseEssai"
$ws.Range("I3").Value = "#Resultat"
$ws.Range("J3").Value = "#NuméroLotReactif"
